$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the two existing result rows down by one: row 3 (16/6 Mon) <- row 2 (15/6 Sun)
# data currently sits in rows 2 and 3; push them down to make room for the new
# 18/6 (Wed) entry at row 2. (.Text is used for the read side since .Value
# round-trips unreliably through this shim.)
$ws.Range("A4").Value = $ws.Range("A3").Text
$ws.Range("B4").Value = $ws.Range("B3").Text
$ws.Range("C4").Value = $ws.Range("C3").Text

$ws.Range("A3").Value = $ws.Range("A2").Text
$ws.Range("B3").Value = $ws.Range("B2").Text
$ws.Range("C3").Value = $ws.Range("C2").Text

# Write the new 4D box result for 18/6/2025 (Wed) into row 2.
$ws.Range("A2").Value = "18/6/2025 (Wed)"
$ws.Range("B2").Value = "2 1 2 1`n3 2 4 7`n0 5 6 9`n5 3 3 8"
$ws.Range("C2").Value = "✅ Direct: 9/3416 (0.26%)`n✅ iBet: 9/188 (4.79%)"

# Re-fit row 2's height back to automatic (writing the multi-line values
# otherwise pins a custom row height on it).
$ws.Rows(2).AutoFit()

# Row 7 gains a styled (wrap-text) placeholder cell in column C, matching the
# other rows in the block.
$ws.Range("C7").WrapText = $true

# Extend the table with a new blank styled row 33, mirroring row 32.
$ws.Range("B33").WrapText = $true
